$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column N (14th column). This shifts the
# old N,O columns (Precisione media / Precisione) to O,P, and pushes
# every column/formula reference to their right accordingly (Excel's
# native Insert behaviour keeps shared-formula refs, dimension, cols
# widths, selection anchors, etc. consistent automatically).
$ws.Columns("N:N").Insert()

# Re-label the "Numero elementi corretti" header (now still in M2) as
# "Elementi corretti", and give the freshly inserted column N2 its new
# header "Elementi corretti (%)".
$ws.Range("M2").Value2 = "Elementi corretti"
$ws.Range("N2").Value2 = "Elementi corretti (%)"

# Fill the new column with the "% of corrected elements" formula
# (= elementi corretti / parole estratte) for the data rows that have
# values (rows 10-24), matching the same relative-formula pattern
# already used by the sheet's other ratio columns.
$ws.Range("N10:N24").Formula = "=M10/D10"

# Match the percentage number format used by the neighbouring
# "Precisione" columns so the new column reuses cell style s="5".
$ws.Range("N10:N24").NumberFormat = "0.00%"

# Restore the original column widths for M (unchanged data, but Excel
# renders it a bit narrower after the insert) and set the new N
# column's width; values follow the workbook's existing width scheme.
$ws.Columns("M:M").ColumnWidth = 18.1640625 - 0.8333333333333
$ws.Columns("N:N").ColumnWidth = 19.6640625 - 0.8333333333333

# Restore the selection to the cell the author ended up on.
$ws.Range("D27").Select()
